# make get_quality smarter. it will assign ambiguous quality given one input
# Log four more FlaskMegaTut time-tracking entries at the bottom of the sheet,
# matching the formatting already used by the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the last existing data row (46) as the formatting template so the new
# rows reuse the same cell styles instead of creating new ones.
$templateA = $ws.Cells.Item(46, 1)
$templateB = $ws.Cells.Item(46, 2)
$templateC = $ws.Cells.Item(46, 3)

$rows = @(47, 48, 49, 50)
$dates = @(43351, 43351, 43354, 43355)
$durations = @(0.048078703703703707, 0.045231481481481484, 0.022210648148148149, 0.025381944444444443)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]

    $cellA = $ws.Cells.Item($r, 1)
    $templateA.Copy()
    $cellA.PasteSpecial(-4122)
    $cellA.Value = "FlaskMegaTut"

    $cellB = $ws.Cells.Item($r, 2)
    $templateB.Copy()
    $cellB.PasteSpecial(-4122)
    $cellB.Value = $dates[$i]

    $cellC = $ws.Cells.Item($r, 3)
    $templateC.Copy()
    $cellC.PasteSpecial(-4122)
    $cellC.Value = $durations[$i]
}

# Reflect where the user was last working in the sheet.
$ws.Range("D28").Select()
